# Adapt column header formatting to respective input file names (#7)
#  - Segmentname_old/_new (etc.) -> Segmentname_FV2310/_FV2404 (etc.)
#  - Wrap the header+data range in an Excel Table ("Table1")
#  - Freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Rename the 21 header cells (A1:U1) from the "_old"/"_new" suffix
#    convention to the "_FV2310"/"_FV2404" format-version suffix.
$headers = @(
  "Segmentname_FV2310",
  "Segmentgruppe_FV2310",
  "Segment_FV2310",
  "Datenelement_FV2310",
  "Segment ID_FV2310",
  "Code_FV2310",
  "Qualifier_FV2310",
  "Beschreibung_FV2310",
  "Bedingungsausdruck_FV2310",
  "Bedingung_FV2310",
  "diff",
  "Segmentname_FV2404",
  "Segmentgruppe_FV2404",
  "Segment_FV2404",
  "Datenelement_FV2404",
  "Segment ID_FV2404",
  "Code_FV2404",
  "Qualifier_FV2404",
  "Beschreibung_FV2404",
  "Bedingungsausdruck_FV2404",
  "Bedingung_FV2404"
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# 2) Turn the used range into a native Excel Table ("Table1") so the
#    header row becomes filterable/structured like the source export.
$usedRange = $ws.Range("A1:U89")
$tbl = $ws.ListObjects.Add(
  [Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange,
  $usedRange,
  $null,
  [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes
)
$tbl.Name = "Table1"

# 3) Freeze the header row (split below row 1) like the refreshed export.
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true

Write-Host "Headers renamed, Table1 created, header row frozen."
